# MAJ staff et feuille de route
# - Stéphane Larose
# - Ajout rencontre production TV / moto le mardi AM

$wb = $excel.ActiveWorkbook

# --- MAR (Tuesday) sheet: insert new meeting row before the current row 9 ---
$wsMar = $wb.Worksheets.Item("MAR")

# Insert a new blank row at position 9 (pushes the former row 9.. down by one)
$wsMar.Rows.Item(9).Insert()

# Fill in the new row's content.
# Note: value assignment order matters because new text gets appended to the
# shared string table in the order the values are first written, and we want
# to reproduce the target shared-string ordering (time, then FR, then EN).
$wsMar.Range("A9").Value = "reunion"
$wsMar.Range("C9").Value = "11:00 - 11:30"
$wsMar.Range("B9").Value = "Réunion production TV <br/>Local TBD"
$wsMar.Range("D9").Value = "TV production meeting <br/>Local TBD"

# The row uses the same "highlighted meeting" row height as similar rows.
$wsMar.Rows.Item(9).RowHeight = 22

# --- Update the active sheet / selection bookkeeping ---
# Previously LUN_AV was the selected tab with B8 selected; the edit moves the
# active tab to MAR with C11 selected.
$wsLunAv = $wb.Worksheets.Item("LUN_AV")
$wsLunAv.Activate() | Out-Null
$wsLunAv.Range("C19").Select() | Out-Null

$wsMar.Activate() | Out-Null
$wsMar.Range("C11").Select() | Out-Null
